# Generate Report for Archive
#
# The status text "Ready for handoff" is updated to "In Translation" on every
# sheet that reports it (Overview!E2:F2, zh-cn!C2, de-de!C2). The Status
# columns are then re-sized to fit the new (shorter) text, matching the
# narrower column width used by the report generator.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# Update the status text wherever it appears.
$overview.Range("E2").Value = "In Translation"
$overview.Range("F2").Value = "In Translation"
$zhcn.Range("C2").Value = "In Translation"
$dede.Range("C2").Value = "In Translation"

# Shrink the Status columns to fit the shorter text (was 17.2159881591797 ->
# now 13.4101845877511 in the report).
$overview.Columns.Item(5).ColumnWidth = 12.5
$overview.Columns.Item(6).ColumnWidth = 12.5
$zhcn.Columns.Item(3).ColumnWidth = 12.5
$dede.Columns.Item(3).ColumnWidth = 12.5
